$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns K and L
$ws.Range("K1").Value = "Relative change (%) Lower UL"
$ws.Range("L1").Value = "Relative change (%) Upper UL"

# New data cells for columns K (Lower UL) and L (Upper UL), rows 2-25
$ws.Range("K2").Value = 0.30000001192092896
$ws.Range("L2").Value = 4.9000000953674317
$ws.Range("K3").Value = -4.5
$ws.Range("L3").Value = 0.80000001192092896
$ws.Range("K4").Value = 5.4000000953674317
$ws.Range("L4").Value = 9.8000001907348633
$ws.Range("K5").Value = -21.899999618530273
$ws.Range("L5").Value = -16
$ws.Range("K6").Value = -27.899999618530273
$ws.Range("L6").Value = -22
$ws.Range("K7").Value = -15.5
$ws.Range("L7").Value = -10
$ws.Range("K8").Value = -11.600000381469727
$ws.Range("L8").Value = -8.6000003814697266
$ws.Range("K9").Value = -16.399999618530273
$ws.Range("L9").Value = -13.399999618530273
$ws.Range("K10").Value = -4.5
$ws.Range("L10").Value = -1.3999999761581421
$ws.Range("K11").Value = -37.099998474121094
$ws.Range("L11").Value = -34.200000762939453
$ws.Range("K12").Value = -41.200000762939453
$ws.Range("L12").Value = -38.299999237060547
$ws.Range("K13").Value = -32.099998474121094
$ws.Range("L13").Value = -29
$ws.Range("K14").Value = 17.100000381469727
$ws.Range("L14").Value = 23.100000381469727
$ws.Range("K15").Value = 19.899999618530273
$ws.Range("L15").Value = 26
$ws.Range("K16").Value = 14.399999618530273
$ws.Range("L16").Value = 21.100000381469727
$ws.Range("K17").Value = 1.8999999761581421
$ws.Range("L17").Value = 6.3000001907348633
$ws.Range("K18").Value = 1.5
$ws.Range("L18").Value = 6.5999999046325684
$ws.Range("K19").Value = 1.5
$ws.Range("L19").Value = 6.4000000953674317
$ws.Range("K20").Value = 5.9000000953674317
$ws.Range("L20").Value = 6.1999998092651367
$ws.Range("K21").Value = 7.8000001907348633
$ws.Range("L21").Value = 7.9000000953674317
$ws.Range("K22").Value = 4.3000001907348633
$ws.Range("L22").Value = 4.5999999046325684
$ws.Range("K23").Value = -56.799999237060547
$ws.Range("L23").Value = -41.200000762939453
$ws.Range("K24").Value = -58.599998474121094
$ws.Range("L24").Value = -43.099998474121094
$ws.Range("K25").Value = -54.400001525878906
$ws.Range("L25").Value = -38.599998474121094
